# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Replaces the "Periodo Mora" (col E) / "Salario Basico" (col F) data block
# (rows 16-50) with the new set of periods (now listed newest -> oldest,
# 2003 down to 1705) and their corresponding updated base-salary values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711","1710","1709","1708","1707","1706","1705"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]

    if ($row -le 34) {
        $ws.Cells.Item($row, 6).Value = 31249
    } else {
        $ws.Cells.Item($row, 6).Value = 29509
    }
}
